# REVER_DailyTracker_NIRMAL.xlsx - "Add files via upload"
# Adds three new daily-tracker rows (21, 22, 23 => sheet rows 22-24) to the
# NOV-2020 sheet: two "Week off" days followed by a completed task day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "NOV-2020"

# --- Row 22: Week off (21 Nov 2020) ------------------------------------
# Style pattern matches existing "Week off" rows (e.g. row 8/9/15/16).
$ws.Range("A8:G8").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 44156
$ws.Cells.Item(22, 4).Value = "Week off"

# --- Row 23: Week off (22 Nov 2020) ------------------------------------
$ws.Range("A8:G8").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = 44157
$ws.Cells.Item(23, 4).Value = "Week off"

# --- Row 24: Completed task (23 Nov 2020) ------------------------------
# Style pattern matches the preceding completed-task rows (e.g. row 21).
$ws.Range("A21:G21").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = 44158
$ws.Cells.Item(24, 3).Value = "B2C/B2B, Sonia and nMVAR"
$ws.Cells.Item(24, 4).Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing & Retesting on B2C/B2B app, Sonia Application (ivc ) and nMVAR"
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = "Completed"
$ws.Rows.Item(24).RowHeight = 45

# --- Selection matches the saved workbook state ------------------------
$ws.Range("E24").Select()

Write-Output "rows 22-24 added to NOV-2020"
